# Generate Report for Handoff
# Updates the "Status" / "In Translation" -> "Ready for handoff" strings and
# their associated "Latest ... Datetime" timestamps, and widens the
# status-adjacent columns on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# E2 (zh-cn status) and F2 (de-de status)
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
# G2 Latest HO Xliff Generate Date
$overview.Range("G2").Value = "2016-08-23 10:40:36"

# Widen columns E and F (status columns) on the Overview sheet
$overview.Columns.Item(5).ColumnWidth = 16.333333333333336
$overview.Columns.Item(6).ColumnWidth = 16.333333333333336

# --- zh-cn sheet ------------------------------------------------------
# C2 Status
$zhcn.Range("C2").Value = "Ready for handoff"
# H2 Latest Handoff Datetime
$zhcn.Range("H2").Value = "2016-08-23 10:40:31"

# Widen column C (Status) on the zh-cn sheet
$zhcn.Columns.Item(3).ColumnWidth = 16.333333333333336

# --- de-de sheet ------------------------------------------------------
# C2 Status
$dede.Range("C2").Value = "Ready for handoff"
# H2 Latest Handoff Datetime
$dede.Range("H2").Value = "2016-08-23 10:40:36"

# Widen column C (Status) on the de-de sheet
$dede.Columns.Item(3).ColumnWidth = 16.333333333333336
